$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "Buy"
$ws.Range("J2").Value = "NIFTY2242117250CE"

$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = "10"

$ws.Range("N2").NumberFormat = "@"
$ws.Range("N2").Value = "2"

$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "15"

$ws.Range("Q2").NumberFormat = "@"
$ws.Range("Q2").Value = "3"
